# Automatische test-sync: 2025-06-20 12:30:50
$wb = $excel.ActiveWorkbook

# "Logs" sheet - append new row 11 with the new mail log entry
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(11, 1).Value = "Uitnodiging voor netwerkevent"
$logs.Cells.Item(11, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(11, 3).Value = "Graag nodig ik u uit voor ons zakelijke netwerkevent volgende maand."
$logs.Cells.Item(11, 4).Value = "Samenwerking / Partnerverzoek"
$logs.Cells.Item(11, 6).Value = "2025-06-20 12:30:28"
$logs.Cells.Item(11, 7).Value = "Nee"

# "Dashboard" sheet - bump the count for "Samenwerking / Partnerverzoek" from 3 to 4
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(2, 2).Value = 4

# Extend the conditional formatting ranges on the "Logs" sheet to cover the new row
$dFormatConditions = $logs.Range("D2:D10").FormatConditions
$dFormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D11"))

$gFormatConditions = $logs.Range("G2:G10").FormatConditions
$gFormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G11"))
